$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to Text format temporarily so numeric-looking strings
# (prices and percentages) are preserved as text, matching the source
# workbook's inlineStr cell type instead of being auto-converted to numbers.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '292.63'
$ws.Range('E2').Value = '-0.63%'
$ws.Range('D3').Value = '40.37'
$ws.Range('E3').Value = '0.72%'
$ws.Range('D4').Value = '5.006'
$ws.Range('E4').Value = '-0.36%'
$ws.Range('E5').Value = '-0.72%'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').Value = '4.289'
$ws.Range('E6').Value = '-0.56%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = '1.563'
$ws.Range('E7').Value = '1.71%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '0.9233'
$ws.Range('E8').Value = '0.09%'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D9').Value = '2.381'
$ws.Range('E9').Value = '-0.76%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.1180'
$ws.Range('E10').Value = '0.03%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1811'
$ws.Range('E11').Value = '2.46%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.04391'
$ws.Range('E12').Value = '5.40%'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').Value = '0.08829'
$ws.Range('E13').Value = '2.17%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.1053'
$ws.Range('E14').Value = '-0.12%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001263'
$ws.Range('E15').Value = '-0.60%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005848'
$ws.Range('E16').Value = '1.18%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.340'
$ws.Range('E17').Value = '-1.03%'
$ws.Range('D18').Value = '0.3327'
$ws.Range('D19').Value = '7.868'
$ws.Range('E19').Value = '3.36%'
$ws.Range('E20').Value = '2.56%'
$ws.Range('D21').Value = '0.2843'
$ws.Range('E21').Value = '1.14%'
$ws.Range('D22').Value = '0.03917'
$ws.Range('E22').Value = '2.39%'
$ws.Range('D23').Value = '0.001260'
$ws.Range('E23').Value = '-1.82%'
$ws.Range('D24').Value = '0.003703'
$ws.Range('E24').Value = '-5.03%'
$ws.Range('D25').Value = '0.0001252'
$ws.Range('E25').Value = '-3.32%'
$ws.Range('D26').Value = '0.0003724'
$ws.Range('E26').Value = '-0.33%'
$ws.Range('D38').Value = '0.02342'
$ws.Range('E38').Value = '1.48%'
$ws.Range('D39').Value = '0.05082'
$ws.Range('E39').Value = '1.79%'
$ws.Range('D40').Value = '0.006018'
$ws.Range('E40').Value = '49.84%'
$ws.Range('D41').Value = '0.007801'
$ws.Range('E41').Value = '1.14%'
$ws.Range('E42').Value = '1.16%'
$ws.Range('D43').Value = '0.007392'
$ws.Range('E43').Value = '-0.12%'
$ws.Range('D44').Value = '0.008060'
$ws.Range('E44').Value = '15.53%'
$ws.Range('D45').Value = '0.2908'
$ws.Range('E45').Value = '-8.68%'
$ws.Range('D46').Value = '0.00006206'
$ws.Range('E46').Value = '-3.59%'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').Value = '-0.28%'
$ws.Range('D48').Value = '0.04727'
$ws.Range('E48').Value = '-81.23%'
$ws.Range('D49').Value = '0.004203'
$ws.Range('E49').Value = '-0.32%'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('E50').Value = '-0.28%'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('E51').Value = '-0.28%'

# Remove the temporary text formatting so no stray style is left behind.
$numRange.ClearFormats()
